$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its text formatting so trailing zeros
# and multi-dot values (e.g. "51.730.51") are preserved, matching the
# original inline-string cell content instead of being auto-converted
# to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '51.730.51'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.952.01'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '351.37'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').Value = '106.26'
$ws.Range('E6').Value = '  -5.42%  '
$ws.Range('D7').Value = '0.555'
$ws.Range('E7').Value = '  -0.72%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.606'
$ws.Range('E9').Value = '  -2.33%  '
$ws.Range('D10').Value = '37.86'
$ws.Range('E10').Value = '  -3.88%  '
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').Value = '0.0849'
$ws.Range('E12').Value = '  -3.80%  '
$ws.Range('D13').Value = '18.98'
$ws.Range('E13').Value = '  -5.78%  '
$ws.Range('D14').Value = '3.399.43'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').Value = '7.60'
$ws.Range('E15').Value = '  -2.26%  '
$ws.Range('D16').Value = '2.940.29'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '0.976'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '51.640.43'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').Value = '3.38'
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').Value = '7.40'
$ws.Range('E20').Value = '  -2.77%  '
$ws.Range('D21').Value = '13.44'
$ws.Range('E21').Value = '  -5.51%  '
$ws.Range('D22').Value = '0.0₃0966'
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('D23').Value = '69.06'
$ws.Range('E23').Value = '  -2.98%  '
$ws.Range('D24').Value = '261.44'
$ws.Range('E24').Value = '  -2.56%  '
$ws.Range('E25').Value = '  -2.65%  '
$ws.Range('E26').Value = '  -3.67%  '
$ws.Range('D27').Value = '26.53'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = '7.31'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('D31').Value = '10.25'
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').Value = '6.11'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = '35.76'
$ws.Range('E33').Value = '  -3.88%  '
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').Value = '2.15'
$ws.Range('E34').Value = '  -4.47%  '
$ws.Range('D35').Value = '50.51'
$ws.Range('E35').Value = '  -4.56%  '
$ws.Range('D36').Value = '0.0430'
$ws.Range('E36').Value = '  -5.17%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').Value = '3.15'
$ws.Range('E38').Value = '  -5.56%  '
$ws.Range('D39').Value = '17.70'
$ws.Range('E39').Value = '  -5.30%  '
$ws.Range('D40').Value = '1.95'
$ws.Range('E40').Value = '  -4.99%  '
$ws.Range('D41').Value = '2.66'
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('E42').Value = '  -1.72%  '
$ws.Range('D43').Value = '123.50'
$ws.Range('E43').Value = '  +11.31%  '
$ws.Range('D44').Value = '22.26'
$ws.Range('E44').Value = '  -3.78%  '
$ws.Range('D45').Value = '2.11'
$ws.Range('E45').Value = '  -3.62%  '
$ws.Range('D46').Value = '2.104.47'
$ws.Range('E46').Value = '  -3.20%  '
$ws.Range('D47').Value = '3.31'
$ws.Range('E47').Value = '  -5.62%  '
$ws.Range('D48').Value = '2.31'
$ws.Range('E48').Value = '  -9.11%  '
$ws.Range('D49').Value = '0.236'
$ws.Range('E49').Value = '  -5.33%  '
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('D51').Value = '0.912'
$ws.Range('E51').Value = '  -3.48%  '
